# Apply the edits described by the commit:
#   "Created vulnerability images and updated README"
#
# Semantic changes made to the template:
#   1. The title paragraph ("Template") no longer uses the Title style;
#      it reverts to the document's default (Normal) style.
#   2. The Normal style now carries its own "space after" of 12pt (240 twips)
#      instead of relying on the document-default spacing.
#   3. The Body Text style drops its own explicit before/after/line spacing
#      override so it again inherits straight from Normal (0pt before,
#      12pt after, single line spacing).
#   4. The Author style's line spacing changes from double (480) to
#      1.5 lines (360).

$d = $word.ActiveDocument

# 1. First paragraph ("Template") loses the Title paragraph style.
$d.Paragraphs(1).Style = "Normal"

# 2. Normal style gains an explicit "space after" of 12pt (240 twips).
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.SpaceAfter = 12

# 3. Body Text reverts to the inherited spacing (no local override):
#    before = 0pt, after = 12pt (matches the new Normal default),
#    single line spacing.
$bodyText = $d.Styles("BodyText")
$bodyText.ParagraphFormat.SpaceBefore = 0
$bodyText.ParagraphFormat.SpaceAfter = 12
$bodyText.ParagraphFormat.LineSpacingRule = 0

# 4. Author style switches from double spacing to 1.5-line spacing.
$author = $d.Styles("Author")
$author.ParagraphFormat.LineSpacingRule = 1
